$d = $word.ActiveDocument

# 1. Remove the "_GoBack" bookmark that currently sits in the title
#    paragraph (between "PDG " and "- Journal de bord").
$d.Bookmarks.Item("_GoBack").Delete()

# 2. Replace the trailing empty paragraph with two new paragraphs:
#      - "Semaine 4 - 05.10.2015" (Titre1 style)
#      - "Rendu du cahier des charges" (Normal style), carrying the
#        "_GoBack" bookmark, collapsed, at the very end of its text.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Text = "Semaine 4 – 05.10.2015"
$lastPara.Style = "Titre1"
$lastPara.Range.InsertParagraphAfter()

$newLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$newLast.Style = "Normal"
# Type a trailing placeholder character so the desired bookmark position
# ("right after the real text") is never the literal last offset of a
# paragraph range, then add the collapsed bookmark there, then strip the
# placeholder back out with a narrow Find/Replace (leaves the already
# placed bookmark untouched).
$newLast.Range.Text = "Rendu du cahier des charges~"

$bmPos = $newLast.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$cleanupRange = $newLast.Range.Duplicate
[void]$cleanupRange.Find.Execute("~", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
